# inicio do backend para criar o fluxo de caixa
#
# Applies to the Kanban workbook:
#  1. "Dados do Roteiro" sheet: mark "View Gerar Relatório" as Concluído,
#     drop the "Testes Automátizados" milestone entirely, and register
#     three new milestones (Desconto nas vendas / Fluxo de caixa / View de
#     relatório de produtos por cliente) at the bottom of the table.
#  2. "Erros" sheet: mark the "Botão voltar" bug (row 10) as Resolvido
#     instead of Identificado.
#  3. Leave the workbook focused on the "Dados do Roteiro" sheet.

$wb = $excel.ActiveWorkbook

$wsRoteiro = $wb.Worksheets.Item("Dados do Roteiro")
$wsMetodologia = $wb.Worksheets.Item("Metodologia Ágil")
$wsErros = $wb.Worksheets.Item("Erros")

# --- 1. Dados do Roteiro --------------------------------------------------

$lo = $wsRoteiro.ListObjects.Item(1)

# "View Gerar Relatório" (row 21) moves from "A Fazer" to "Concluído".
$wsRoteiro.Range("B21").Value = "Concluído"

# Drop the "Testes Automátizados" row (row 22) completely - this shifts
# every following row (and the footer note below the table) up by one.
$wsRoteiro.Rows.Item(22).Delete()

# Append three brand-new milestones at the bottom of the table. Insert the
# rows explicitly (rather than relying on ListRows.Add, which does not push
# the trailing footer row down) so the merged footer note keeps sliding
# down underneath the growing table.
$wsRoteiro.Rows.Item(31).Insert()
$wsRoteiro.Rows.Item(31).Insert()
$wsRoteiro.Rows.Item(31).Insert()
$lo.Resize($wsRoteiro.Range("B2:C34"))

$wsRoteiro.Range("B30").Value = "A Fazer"
$wsRoteiro.Range("C30").Value = "Desconto nas vendas"

$wsRoteiro.Range("B31").Value = "Em Progresso"
$wsRoteiro.Range("C31").Value = "Fluxo de caixa"

$wsRoteiro.Range("B32").Value = "Concluído"
$wsRoteiro.Range("C32").Value = "View de relatório de produtos por cliente"

# --- 2. Erros --------------------------------------------------------------

# The "O botão de voltar..." bug (row 10) has been fixed.
$wsErros.Range("B10").Value = "Resolvido"

# --- 3. View / selection bookkeeping ---------------------------------------

$wsMetodologia.Activate() | Out-Null
$wsMetodologia.Range("C3").Select() | Out-Null

$wsErros.Activate() | Out-Null
$wsErros.Range("C14").Select() | Out-Null

$wsRoteiro.Activate() | Out-Null
$wsRoteiro.Range("C22").Select() | Out-Null
